# Updates the cryptos price/volume table (rows 2-51) to the refreshed values
# pulled by the GitHub Actions job on Fri Jan 19 10:54:34 UTC 2024.
# Row 43/44 additionally swap places (EnergySwap now ranks above Maker).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.458.09'
$ws.Range("E2").Value = '  -2.51%  '
$ws.Range("D3").Value = '2.487.59'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.04'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.79'
$ws.Range("E6").Value = '  -4.05%  '
$ws.Range("E7").Value = '  -2.62%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -3.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.63'
$ws.Range("E10").Value = '  -4.86%  '
$ws.Range("E11").Value = '  -2.53%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("E13").Value = '  -3.08%  '
$ws.Range("D14").Value = '2.867.30'
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.46'
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '2.457.35'
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").Value = '41.417.22'
$ws.Range("E18").Value = '  -2.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.34'
$ws.Range("E19").Value = '  -3.85%  '
$ws.Range("D20").Value = '0.0₃0925'
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.31'
$ws.Range("E21").Value = '  -7.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.01'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.41'
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("E24").Value = '  -3.36%  '
$ws.Range("E25").Value = '  -4.32%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.24'
$ws.Range("E27").Value = '  -4.75%  '
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.77'
$ws.Range("E29").Value = '  -2.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.56'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.00'
$ws.Range("E31").Value = '  -3.15%  '
$ws.Range("E32").Value = '  -5.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.56'
$ws.Range("E33").Value = '  -3.89%  '
$ws.Range("E34").Value = '  -3.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.02'
$ws.Range("E35").Value = '  +2.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0753'
$ws.Range("E36").Value = '  -4.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.10'
$ws.Range("E37").Value = '  -1.48%  '
$ws.Range("E38").Value = '  -3.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.115'
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("E41").Value = '  +1.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.63'
$ws.Range("E43").Value = '  -10.05%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.005.35'
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("E45").Value = '  -3.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.04'
$ws.Range("E46").Value = '  -7.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.84'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("D48").Value = '2.731.09'
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '70.02'
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '97.30'
$ws.Range("E50").Value = '  -3.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.76'
$ws.Range("E51").Value = '  -5.06%  '
